$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "[Victor-Usin. CNC, Leandro-M.S.R. ar Cond., Leandro-Mec. Manut.Equip. ind., Ludoff-Coman. Hidraulicos]"
$ws.Range("C2").Value = "[Elcio Dec.-Cont.Lóg.Prog CLP, Joel L.-Tec. Fundição, Rogério-Retífica, -]"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"

# Row 3
$ws.Range("B3").Value = "[Ivan-Trat. Termicos, Valmir-Calderaria, Euclides-Tec. Soldagem, Aderci-Fresagem]"
$ws.Range("C3").Value = "[Aselmo-M. Motor Endot., Gisele-Ens. Dest. Não Desti., Ismail-Metrologia 2, Paulo Rob.-M.A.Comp.CAD / CAM]"
$ws.Range("D3").Value = "[Ludoff-Coman. Hidraulicos, Victor-Usin. CNC, Joel L.-Tec. Fundição, Rogério-Retífica]"
$ws.Range("E3").Value = "-"

# Row 4
$ws.Range("B4").Value = "[Ivan-Trat. Termicos, Valmir-Calderaria, Euclides-Tec. Soldagem, Aderci-Fresagem]"
$ws.Range("C4").Value = "[Aselmo-M. Motor Endot., Gisele-Ens. Dest. Não Desti., Ismail-Metrologia 2, Paulo Rob.-M.A.Comp.CAD / CAM]"
$ws.Range("D4").Value = "[Ludoff-Coman. Hidraulicos, Victor-Usin. CNC, Joel L.-Tec. Fundição, Rogério-Retífica]"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "Nilton-Elem. Máqui"

# Row 6
$ws.Range("B6").Value = "[Ivan-Trat. Termicos, Valmir-Calderaria, Euclides-Tec. Soldagem, Aderci-Fresagem]"
$ws.Range("C6").Value = "[Aselmo-M. Motor Endot., Gisele-Ens. Dest. Não Desti., Ismail-Metrologia 2, Paulo Rob.-M.A.Comp.CAD / CAM]"
$ws.Range("D6").Value = "[Elcio Dec.-Cont.Lóg.Prog CLP, -, Joel L.-Tec. Fundição, Rogério-Retífica]"
$ws.Range("E6").Value = "-"

# Row 7
$ws.Range("B7").Value = "[Ivan-Trat. Termicos, Valmir-Calderaria, Euclides-Tec. Soldagem, Aderci-Fresagem]"
$ws.Range("C7").Value = "[Leandro-Mec. Manut.Equip. ind., Leandro-M.S.R. ar Cond., Ismail-Metrologia 2, Paulo Rob.-M.A.Comp.CAD / CAM]"
$ws.Range("D7").Value = "[Elcio Dec.-Cont.Lóg.Prog CLP, -, -, -]"
$ws.Range("E7").Value = "[Mayra-Tec. Mat. Não Metal., Mayra-Tec. Mat. Não Metal.]"

# Row 8
$ws.Range("B8").Value = "[Gisele-Ens. Dest. Não Desti., Aselmo-M. Motor Endot., Leandro-M.S.R. ar Cond., Leandro-Mec. Manut.Equip. ind.]"
$ws.Range("C8").Value = "[Leandro-Mec. Manut.Equip. ind., Leandro-M.S.R. ar Cond., Ludoff-Coman. Hidraulicos, Victor-Usin. CNC]"
$ws.Range("D8").Value = "[Elcio Dec.-Cont.Lóg.Prog CLP, -, -, -]"
$ws.Range("E8").Value = "-"
